$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.387.63'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').Value = '3.909.97'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.718'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.173'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.19%  '
$ws.Range('E11').Value = '  -5.67%  '
$ws.Range('E12').Value = '  -2.24%  '
$ws.Range('D13').Value = '4.537.84'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('D15').Value = '3.919.82'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.23'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.99'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.44%  '
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = '69.351.78'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '429.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  -3.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '88.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.82%  '
$ws.Range('E27').Value = '  -3.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '678.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.33%  '
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '69.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.439'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.73%  '
$ws.Range('D35').Value = '0.0₃0868'
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.95'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '40.03'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.149'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.34%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0481'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.10%  '
$ws.Range('E44').Value = '  -6.16%  '
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').Value = '0.0₆0361'
$ws.Range('E47').Value = '  +13.90%  '
$ws.Range('E48').Value = '  +7.15%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000275'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +13.61%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.750.64'
$ws.Range('E50').Value = '  +13.89%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.07%  '
